$wb = $excel.ActiveWorkbook

# --- Master sheet: fix region split for the Gsteel unit (EU27 / RoW) ---
$wsMaster = $wb.Worksheets.Item("Master")

# Row 2 was tagged "GLOBAL" - it should instead be the EU27 split of Gsteel.
$wsMaster.Range("A2").Value = "EU27"

# Add the missing RoW row for the Gsteel unit (this was being dropped as a
# "duplicate" before the fix).
$wsMaster.Range("A3").Value = "RoW"
$wsMaster.Range("B3").Value = "Green steelmaking"
$wsMaster.Range("C3").Value = "Green steel"
$wsMaster.Range("D3").Value = "Gsteel RoW"
$wsMaster.Range("E3").Value = 1
$wsMaster.Range("F3").Value = "Mton"
$wsMaster.Range("G3").Value = 1

# Match row 2's number formatting on the new FU quantity cell.
$wsMaster.Range("E2").Copy()
$wsMaster.Range("E3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsMaster.Activate()
$wsMaster.Range("J3").Select()

# --- Regions Map: update the remembered selection ---
$wsRegions = $wb.Worksheets.Item("Regions Map")
$wsRegions.Activate()
$wsRegions.Range("A2").Select()

# --- Gsteel: update the remembered selection ---
$wsGsteel = $wb.Worksheets.Item("Gsteel")
$wsGsteel.Activate()
$wsGsteel.Range("D4").Select()

# Leave the Master tab as the active / selected tab.
$wsMaster.Activate()
